$d = $word.ActiveDocument

$d.Content.Find.Execute("810÷7=115, 5", $true, $false, $false, $false, $false, $true, 1, $false, "393÷9=43, 6", 2) | Out-Null
$d.Content.Find.Execute("447÷8=55, 7", $true, $false, $false, $false, $false, $true, 1, $false, "790÷9=87, 7", 2) | Out-Null
$d.Content.Find.Execute("713÷7=101, 6", $true, $false, $false, $false, $false, $true, 1, $false, "194÷9=21, 5", 2) | Out-Null
$d.Content.Find.Execute("874÷5=174, 4", $true, $false, $false, $false, $false, $true, 1, $false, "692÷6=115, 2", 2) | Out-Null
$d.Content.Find.Execute("101÷5=20, 1", $true, $false, $false, $false, $false, $true, 1, $false, "431÷4=107, 3", 2) | Out-Null
$d.Content.Find.Execute("621÷7=88, 5", $true, $false, $false, $false, $false, $true, 1, $false, "564÷7=80, 4", 2) | Out-Null
$d.Content.Find.Execute("622÷5=124, 2", $true, $false, $false, $false, $false, $true, 1, $false, "352÷7=50, 2", 2) | Out-Null
$d.Content.Find.Execute("663÷7=94, 5", $true, $false, $false, $false, $false, $true, 1, $false, "861÷5=172, 1", 2) | Out-Null
$d.Content.Find.Execute("482÷4=120, 2", $true, $false, $false, $false, $false, $true, 1, $false, "812÷6=135, 2", 2) | Out-Null
$d.Content.Find.Execute("915÷2=457, 1", $true, $false, $false, $false, $false, $true, 1, $false, "216÷9=24, 0", 2) | Out-Null
$d.Content.Find.Execute("482÷3=160, 2", $true, $false, $false, $false, $false, $true, 1, $false, "465÷3=155, 0", 2) | Out-Null
$d.Content.Find.Execute("528÷6=88, 0", $true, $false, $false, $false, $false, $true, 1, $false, "932÷7=133, 1", 2) | Out-Null
$d.Content.Find.Execute("331÷9=36, 7", $true, $false, $false, $false, $false, $true, 1, $false, "985÷2=492, 1", 2) | Out-Null
$d.Content.Find.Execute("647÷8=80, 7", $true, $false, $false, $false, $false, $true, 1, $false, "645÷3=215, 0", 2) | Out-Null
$d.Content.Find.Execute("981÷4=245, 1", $true, $false, $false, $false, $false, $true, 1, $false, "470÷5=94, 0", 2) | Out-Null
$d.Content.Find.Execute("880÷5=176, 0", $true, $false, $false, $false, $false, $true, 1, $false, "783÷6=130, 3", 2) | Out-Null
$d.Content.Find.Execute("722÷3=240, 2", $true, $false, $false, $false, $false, $true, 1, $false, "203÷7=29, 0", 2) | Out-Null
$d.Content.Find.Execute("519÷8=64, 7", $true, $false, $false, $false, $false, $true, 1, $false, "214÷9=23, 7", 2) | Out-Null
$d.Content.Find.Execute("477÷7=68, 1", $true, $false, $false, $false, $false, $true, 1, $false, "253÷9=28, 1", 2) | Out-Null
$d.Content.Find.Execute("632÷5=126, 2", $true, $false, $false, $false, $false, $true, 1, $false, "622÷4=155, 2", 2) | Out-Null
$d.Content.Find.Execute("977÷8=122, 1", $true, $false, $false, $false, $false, $true, 1, $false, "817÷3=272, 1", 2) | Out-Null
$d.Content.Find.Execute("984÷4=246, 0", $true, $false, $false, $false, $false, $true, 1, $false, "675÷2=337, 1", 2) | Out-Null
$d.Content.Find.Execute("459÷6=76, 3", $true, $false, $false, $false, $false, $true, 1, $false, "531÷4=132, 3", 2) | Out-Null
$d.Content.Find.Execute("480÷9=53, 3", $true, $false, $false, $false, $false, $true, 1, $false, "622÷8=77, 6", 2) | Out-Null
$d.Content.Find.Execute("237÷3=79, 0", $true, $false, $false, $false, $false, $true, 1, $false, "253÷3=84, 1", 2) | Out-Null
